$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.824.23'
$ws.Range("E2").Value = '  +3.41%  '
$ws.Range("D3").Value = '2.264.96'
$ws.Range("E3").Value = '  -0.25%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.78'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.25'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.47%  '
$ws.Range("E7").Value = '  -1.25%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.31'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +3.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0814'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.08'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("E13").Value = '  -1.15%  '
$ws.Range("D14").Value = '2.610.38'
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = '2.262.85'
$ws.Range("E15").Value = '  -0.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.60'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").Value = '46.837.32'
$ws.Range("E17").Value = '  +3.75%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.790'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.76'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.09%  '
$ws.Range("D20").Value = '0.0₃0954'
$ws.Range("E20").Value = '  +4.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.81'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.35'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '246.88'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.87%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.79'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '41.72'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.25'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.63'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.50'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +4.43%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.84'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +10.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '145.92'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.29'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +13.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.38'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0767'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.77%  '
$ws.Range("E36").Value = '  +10.36%  '
$ws.Range("E37").Value = '  -2.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.80'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +17.47%  '
$ws.Range("E39").Value = '  -4.17%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.84'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.26%  '
$ws.Range("E41").Value = '  -4.99%  '
$ws.Range("E42").Value = '  -2.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '91.05'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +19.51%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.784.13'
$ws.Range("E45").Value = '  +1.21%  '
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.88'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '71.24'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.75%  '
$ws.Range("E48").Value = '  -3.66%  '
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '94.13'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.81'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.69%  '
